$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (all values in this sheet are
# stored as text, e.g. "41.720.38" or "  -2.71%  ", so we force the
# NumberFormat to Text before assignment to stop Excel from
# reinterpreting them as numbers and dropping formatting like trailing
# zeros or the surrounding spaces).
$updates = @{
    "D2" = '41.720.38'
    "E2" = '  -2.71%  '
    "D3" = '2.282.22'
    "E3" = '  -3.59%  '
    "E4" = '  -0.09%  '
    "D5" = '314.50'
    "E5" = '  -1.21%  '
    "D6" = '102.44'
    "E6" = '  -5.91%  '
    "D7" = '0.626'
    "E7" = '  -1.56%  '
    "E8" = '  +0.03%  '
    "E9" = '  -4.06%  '
    "D10" = '38.65'
    "E10" = '  -8.16%  '
    "D11" = '0.0901'
    "E11" = '  -2.85%  '
    "E12" = '  -4.49%  '
    "D13" = '0.105'
    "E13" = '  -1.02%  '
    "E14" = '  -5.14%  '
    "D15" = '15.23'
    "E15" = '  -5.86%  '
    "D16" = '2.627.46'
    "E16" = '  -3.58%  '
    "D17" = '2.283.17'
    "E17" = '  -2.69%  '
    "D18" = '41.694.81'
    "E18" = '  -2.64%  '
    "D19" = '7.39'
    "E19" = '  -4.32%  '
    "E20" = '  -1.85%  '
    "D21" = '3.61'
    "E21" = '  -2.09%  '
    "D22" = '72.99'
    "E22" = '  -4.35%  '
    "D23" = '279.52'
    "E23" = '  +8.46%  '
    "D24" = '10.12'
    "E24" = '  +6.88%  '
    "E25" = '  -4.10%  '
    "E26" = '  +0.64%  '
    "D27" = '10.69'
    "E27" = '  -7.06%  '
    "E28" = '  +6.19%  '
    "D29" = '22.87'
    "E29" = '  -0.54%  '
    "D30" = '162.63'
    "E30" = '  -5.89%  '
    "D31" = '34.86'
    "E31" = '  -6.96%  '
    "E32" = '  -3.36%  '
    "E33" = '  -2.39%  '
    "E34" = '  -5.07%  '
    "D35" = '0.132'
    "E35" = '  +0.06%  '
    "E36" = '  -5.95%  '
    "E37" = '  -4.19%  '
    "D38" = '2.87'
    "E38" = '  +6.16%  '
    "E39" = '  -5.76%  '
    "D40" = '3.61'
    "E40" = '  -8.65%  '
    "D41" = '99.78'
    "E41" = '  +15.47%  '
    "E42" = '  -4.35%  '
    "D43" = '69.32'
    "E43" = '  -3.60%  '
    "E44" = '  -0.16%  '
    "E45" = '  -7.95%  '
    "D46" = '115.79'
    "E46" = '  +2.34%  '
    "D47" = '11.84'
    "E47" = '  -4.67%  '
    "D48" = '8.90'
    "E48" = '  -3.89%  '
    "D49" = '75.57'
    "E49" = '  -2.04%  '
    "D50" = '5.23'
    "E50" = '  -7.00%  '
    "E51" = '  -4.45%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
